# Apply corrections to statsec / ggw7_type codes and their "eerste wijk" flag
# as described in the commit message ("Correctie van ... statsec namen").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 - niscode 12035
$ws.Range("D41").Value = 2.00
$ws.Range("E41").Value = "12035A0"

# Row 160 - niscode 31022
$ws.Range("D160").Value = 1.00
$ws.Range("E160").Value = "31022BAB"

# Row 164 - niscode 31043
$ws.Range("D164").Value = 1.00
$ws.Range("E164").Value = "31043DUI"

# Row 176 - niscode 33040
$ws.Range("D176").Value = 1.00
$ws.Range("E176").Value = "33040BIKS"

# Row 180 - niscode 34009
$ws.Range("D180").Value = 1.00
$ws.Range("E180").Value = "34009BELG"

# Row 182 - niscode 34022 (Kortrijk) - fix duplicated niscode prefix in the wijkcode
$ws.Range("E182").Value = "34022WA"

# Row 185 - niscode 34027
$ws.Range("D185").Value = 1.00
$ws.Range("E185").Value = "34027BAR"

# Row 186 - niscode 34040
$ws.Range("D186").Value = 1.00
$ws.Range("E186").Value = "34040BEVL"

# Row 188 - niscode 34042
$ws.Range("D188").Value = 1.00
$ws.Range("E188").Value = "34042CENT"

# Row 306 - niscode 72042 (was mistakenly linked to 71047A0)
$ws.Range("D306").Value = 1.00
$ws.Range("E306").Value = "72042EL"

$wb.Save()
